# Add a "Save" column (H) to the s_vals worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - copy style from the other header cells (e.g. G1) then set value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Values for H2:H13 as described in the diff.
$saveValues = @(0, 0, 0, 0, 1, 0, 0, 0, 0, 1, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
